$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = ""
$ws.Range("C1").Value = ""

$ws.Range("A3").Value = 5
$ws.Range("C3").Value = 1

$ws.Range("C7").Select()
